$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.194.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.81%  '

$ws.Range("D3").Value = '''1.850.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("D5").Value = '''0.7028'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.75%  '

$ws.Range("D6").Value = '''239.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.24%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").Value = '''0.3055'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.58%  '

$ws.Range("D9").Value = '''0.07429'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.53%  '

$ws.Range("D10").Value = '''23.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.43%  '

$ws.Range("D11").Value = '''0.08157'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.10%  '

$ws.Range("D12").Value = '''1.883.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").Value = '''0.7286'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.74%  '

$ws.Range("D14").Value = '''5.219'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").Value = '''89.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.27%  '

$ws.Range("D16").Value = '''29.418.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.03%  '

$ws.Range("D17").Value = '''5.787'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.83%  '

$ws.Range("D18").Value = '''238.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.30%  '

$ws.Range("E19").Value = '  -3.16%  '

$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '''1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''7.603'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.54%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '''9.015'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.85%  '

$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = '''0.1457'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.50%  '

$ws.Range("D26").Value = '''160.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.65%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''18.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.00%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''1.979'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''1.411'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.30%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.517'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.493'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.68%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''4.014'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.67%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.05204'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.56%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.188'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.97%  '

$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = '''1.046'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.64%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.7081'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.66%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.663'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.11%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01870'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.63%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.679'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.84%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.9419'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.87%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''6.049'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.4314'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.51%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '''1.068.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.87%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''70.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.33%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''103.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '''2.028.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.750'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.68%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '''7.058'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.86%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.137'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.43%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.05878'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.71%  '

